$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 2
$ws_ALC.Range("H2").Value = 1766.3334
$ws_ALC.Range("I2").Value = 149.5
$ws_ALC.Range("K2").Value = 149.5
$ws_ALC.Range("M2").Value = -36.5

# ALC row 9
$ws_ALC.Range("H9").Value = 142.6
$ws_ALC.Range("J9").Value = 132.66667
$ws_ALC.Range("L9").Value = 132.66667
$ws_ALC.Range("N9").Value = -470.66667

# ALC row 10
$ws_ALC.Range("H10").Value = 3999
$ws_ALC.Range("J10").Value = 4999
$ws_ALC.Range("L10").Value = 4999
$ws_ALC.Range("N10").Value = -5585

# ALC row 33
$ws_ALC.Range("H33").Value = 1425505
$ws_ALC.Range("I33").Value = 2458741.5
$ws_ALC.Range("J33").Value = 4804.875
$ws_ALC.Range("K33").Value = 2458741.5
$ws_ALC.Range("L33").Value = 4804.875
$ws_ALC.Range("M33").Value = -2458512.5
$ws_ALC.Range("N33").Value = -5262.875

# ALC row 82
$ws_ALC.Range("H82").Value = 1304.75
$ws_ALC.Range("J82").Value = 0
$ws_ALC.Range("L82").Value = 0
$ws_ALC.Range("N82").ClearContents()

# ALC row 85
$ws_ALC.Range("H85").Value = 1304.75
$ws_ALC.Range("J85").Value = 0
$ws_ALC.Range("L85").Value = 0
$ws_ALC.Range("N85").ClearContents()

# ALC row 116
$ws_ALC.Range("H116").Value = 4321.8076
$ws_ALC.Range("I116").Value = 2500.875
$ws_ALC.Range("J116").Value = 5131.1113
$ws_ALC.Range("K116").Value = 2500.875
$ws_ALC.Range("L116").Value = 5131.1113
$ws_ALC.Range("M116").Value = 941.125
$ws_ALC.Range("N116").Value = -12015.1113

# ALC row 135
$ws_ALC.Range("H135").Value = 2261.3635
$ws_ALC.Range("I135").Value = 2337.6
$ws_ALC.Range("K135").Value = 21038.4
$ws_ALC.Range("M135").Value = -18503.4

# ALC row 137
$ws_ALC.Range("H137").Value = 9231.550999999999
$ws_ALC.Range("I137").Value = 4600.2383
$ws_ALC.Range("J137").Value = 16435.814
$ws_ALC.Range("K137").Value = 13800.7149
$ws_ALC.Range("L137").Value = 49307.442
$ws_ALC.Range("M137").Value = -11250.7149
$ws_ALC.Range("N137").Value = -54407.442

# ALC row 138
$ws_ALC.Range("H138").Value = 1483.0303
$ws_ALC.Range("I138").Value = 1031.4
$ws_ALC.Range("K138").Value = 3094.2
$ws_ALC.Range("M138").Value = 2045.8

# ARM row 45
$ws_ARM.Range("H45").Value = 10245.823
$ws_ARM.Range("I45").Value = 10636.25
$ws_ARM.Range("K45").Value = 10636.25
$ws_ARM.Range("M45").Value = -10259.25

# ARM row 61
$ws_ARM.Range("H61").Value = 4502.838
$ws_ARM.Range("I61").Value = 3140.7097
$ws_ARM.Range("K61").Value = 3140.7097
$ws_ARM.Range("M61").Value = -2928.7097

# ARM row 74
$ws_ARM.Range("H74").Value = 10809.454
$ws_ARM.Range("I74").Value = 11989.852
$ws_ARM.Range("J74").Value = 5497.6665
$ws_ARM.Range("K74").Value = 11989.852
$ws_ARM.Range("L74").Value = 5497.6665
$ws_ARM.Range("M74").Value = -11115.852
$ws_ARM.Range("N74").Value = -7245.6665

# ARM row 77
$ws_ARM.Range("H77").Value = 10809.454
$ws_ARM.Range("I77").Value = 11989.852
$ws_ARM.Range("J77").Value = 5497.6665
$ws_ARM.Range("K77").Value = 59949.26
$ws_ARM.Range("L77").Value = 27488.3325
$ws_ARM.Range("M77").Value = -55581.26
$ws_ARM.Range("N77").Value = -36224.3325

# ARM row 104
$ws_ARM.Range("H104").Value = 0
$ws_ARM.Range("J104").Value = 0
$ws_ARM.Range("L104").Value = 0
$ws_ARM.Range("N104").ClearContents()

# ARM row 132
$ws_ARM.Range("H132").Value = 3223.0667
$ws_ARM.Range("I132").Value = 2240.756
$ws_ARM.Range("K132").Value = 6722.268
$ws_ARM.Range("M132").Value = -4192.268

# ARM row 136
$ws_ARM.Range("H136").Value = 4502.838
$ws_ARM.Range("I136").Value = 3140.7097
$ws_ARM.Range("K136").Value = 9422.1291
$ws_ARM.Range("M136").Value = -6872.1291

# BSM row 20
$ws_BSM.Range("H20").Value = 16490.75
$ws_BSM.Range("I20").Value = 26114.77
$ws_BSM.Range("J20").Value = 5116.909
$ws_BSM.Range("K20").Value = 26114.77
$ws_BSM.Range("L20").Value = 5116.909
$ws_BSM.Range("M20").Value = -25867.77
$ws_BSM.Range("N20").Value = -5610.909

# BSM row 107
$ws_BSM.Range("H107").Value = 1025
$ws_BSM.Range("I107").Value = 1025
$ws_BSM.Range("K107").Value = 1025
$ws_BSM.Range("M107").Value = 895

# BSM row 134
$ws_BSM.Range("H134").Value = 6139.1406
$ws_BSM.Range("I134").Value = 2969.551
$ws_BSM.Range("K134").Value = 8908.653
$ws_BSM.Range("M134").Value = -6373.653

# CRP row 16
$ws_CRP.Range("H16").Value = 84898.414
$ws_CRP.Range("I16").Value = 1222.8334
$ws_CRP.Range("J16").Value = 168574
$ws_CRP.Range("K16").Value = 1222.8334
$ws_CRP.Range("L16").Value = 168574
$ws_CRP.Range("M16").Value = -935.8334
$ws_CRP.Range("N16").Value = -169148

# CRP row 58
$ws_CRP.Range("H58").Value = 2920.1462
$ws_CRP.Range("I58").Value = 1163.862
$ws_CRP.Range("K58").Value = 1163.862
$ws_CRP.Range("M58").Value = -960.8620000000001

# CRP row 113
$ws_CRP.Range("H113").Value = 84898.414
$ws_CRP.Range("I113").Value = 1222.8334
$ws_CRP.Range("J113").Value = 168574
$ws_CRP.Range("K113").Value = 1222.8334
$ws_CRP.Range("L113").Value = 168574
$ws_CRP.Range("M113").Value = 947.1666
$ws_CRP.Range("N113").Value = -172914

# CRP row 132
$ws_CRP.Range("H132").Value = 25226.838
$ws_CRP.Range("I132").Value = 15966.171
$ws_CRP.Range("J132").Value = 43307.19
$ws_CRP.Range("K132").Value = 47898.513
$ws_CRP.Range("L132").Value = 129921.57
$ws_CRP.Range("M132").Value = -45368.513
$ws_CRP.Range("N132").Value = -134981.57

# CRP row 134
$ws_CRP.Range("H134").Value = 2012
$ws_CRP.Range("I134").Value = 1555.7693
$ws_CRP.Range("K134").Value = 4667.3079
$ws_CRP.Range("M134").Value = -2132.3079

# CRP row 136
$ws_CRP.Range("H136").Value = 2920.1462
$ws_CRP.Range("I136").Value = 1163.862
$ws_CRP.Range("K136").Value = 3491.586
$ws_CRP.Range("M136").Value = -941.5860000000002

# CUL row 5
$ws_CUL.Range("H5").Value = 1655.4584
$ws_CUL.Range("I5").Value = 904.2222
$ws_CUL.Range("J5").Value = 2106.2
$ws_CUL.Range("K5").Value = 2712.6666
$ws_CUL.Range("L5").Value = 6318.599999999999
$ws_CUL.Range("M5").Value = -2600.6666
$ws_CUL.Range("N5").Value = -6542.599999999999

# CUL row 12
$ws_CUL.Range("H12").Value = 9999
$ws_CUL.Range("J12").Value = 9999
$ws_CUL.Range("L12").Value = 29997
$ws_CUL.Range("N12").Value = -30343

# CUL row 28
$ws_CUL.Range("H28").Value = 304.5
$ws_CUL.Range("I28").Value = 304.5
$ws_CUL.Range("K28").Value = 913.5
$ws_CUL.Range("M28").Value = -681.5

# CUL row 46
$ws_CUL.Range("H46").Value = 698.75
$ws_CUL.Range("I46").Value = 598.3333
$ws_CUL.Range("K46").Value = 1794.9999
$ws_CUL.Range("M46").Value = -1703.9999

# CUL row 131
$ws_CUL.Range("H131").Value = 3583.5247
$ws_CUL.Range("I131").Value = 884.625
$ws_CUL.Range("J131").Value = 3990.9058
$ws_CUL.Range("K131").Value = 2653.875
$ws_CUL.Range("L131").Value = 11972.7174
$ws_CUL.Range("M131").Value = 2386.125
$ws_CUL.Range("N131").Value = -22052.7174

# CUL row 132
$ws_CUL.Range("H132").Value = 4763543
$ws_CUL.Range("I132").Value = 2439.7
$ws_CUL.Range("J132").Value = 9091819
$ws_CUL.Range("K132").Value = 21957.3
$ws_CUL.Range("L132").Value = 81826371
$ws_CUL.Range("M132").Value = -19427.3
$ws_CUL.Range("N132").Value = -81831431

# CUL row 135
$ws_CUL.Range("H135").Value = 1655.4584
$ws_CUL.Range("I135").Value = 904.2222
$ws_CUL.Range("J135").Value = 2106.2
$ws_CUL.Range("K135").Value = 8137.999800000001
$ws_CUL.Range("L135").Value = 18955.8
$ws_CUL.Range("M135").Value = -5602.999800000001
$ws_CUL.Range("N135").Value = -24025.8

# GSM row 5
$ws_GSM.Range("H5").Value = 9000
$ws_GSM.Range("I5").Value = 9000
$ws_GSM.Range("K5").Value = 9000
$ws_GSM.Range("M5").Value = -8888

# GSM row 103
$ws_GSM.Range("H103").Value = 0
$ws_GSM.Range("J103").Value = 0
$ws_GSM.Range("L103").Value = 0
$ws_GSM.Range("N103").ClearContents()

# LTW row 136
$ws_LTW.Range("H136").Value = 3748.9033
$ws_LTW.Range("I136").Value = 3068.92
$ws_LTW.Range("K136").Value = 9206.76
$ws_LTW.Range("M136").Value = -6656.76

# WVR row 132
$ws_WVR.Range("H132").Value = 24757.988
$ws_WVR.Range("I132").Value = 24865.88
$ws_WVR.Range("J132").Value = 24538.482
$ws_WVR.Range("K132").Value = 74597.64
$ws_WVR.Range("L132").Value = 73615.446
$ws_WVR.Range("M132").Value = -72067.64
$ws_WVR.Range("N132").Value = -78675.446

# WVR row 136
$ws_WVR.Range("H136").Value = 2820.2
$ws_WVR.Range("I136").Value = 1336.7
$ws_WVR.Range("K136").Value = 4010.1
$ws_WVR.Range("M136").Value = -1460.1
